$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phaseshift_metadata")

# --- New data rows (10, 11, 12) -------------------------------------------
# Row 10: 2024-04-25_D_e.dat
$ws.Range("A10").Value = "2024-04-25_D_e.dat"
$ws.Range("B10").Value = "D"
$ws.Range("C10").Value = 0
$ws.Range("D10").Value = "2024-04-25"
$ws.Range("E10").Value = 2
$ws.Range("F10").Value = 1.8
$ws.Range("G10").Value = 50
$ws.Range("H10").Value = 202.14

# Row 11: 2024-05-02_B_e.dat
$ws.Range("A11").Value = "2024-05-02_B_e.dat"
$ws.Range("B11").Value = "B"
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = "2024-05-02"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = 1.8
$ws.Range("G11").Value = 20
$ws.Range("H11").Value = 202

# Row 12: 2024-05-10_V_e.dat
$ws.Range("A12").Value = "2024-05-10_V_e.dat"
$ws.Range("B12").Value = "V"
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = "2024-05-10"
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 1.8
$ws.Range("G12").Value = 20
$ws.Range("H12").Value = 202

# Match the date formatting already used by column D (m/d/yy, style s="1")
$ws.Range("D10:D12").NumberFormat = "m/d/yy"

# --- Column width (col A) --------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.17

# --- Selection --------------------------------------------------------------
$ws.Range("E21").Select()
